$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.221.78'
$ws.Range("E2").Value = '  +3.39%  '
$ws.Range("D3").Value = '1.813.84'
$ws.Range("E3").Value = '  +4.87%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '329.63'
$ws.Range("E5").Value = '  +2.45%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4450'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.28%  '
$ws.Range("D8").Value = '0.3705'
$ws.Range("E8").Value = '  +3.57%  '
$ws.Range("D9").Value = '44.67'
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("D10").Value = '0.07707'
$ws.Range("E10").Value = '  +4.06%  '
$ws.Range("D11").Value = '1.128'
$ws.Range("E11").Value = '  +2.04%  '
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").Value = '22.08'
$ws.Range("E13").Value = '  +3.51%  '
$ws.Range("E14").Value = '  +4.00%  '
$ws.Range("D15").Value = '7.559'
$ws.Range("E15").Value = '  +6.59%  '
$ws.Range("D16").Value = '1.832.93'
$ws.Range("E16").Value = '  +6.00%  '
$ws.Range("D17").Value = '92.81'
$ws.Range("E17").Value = '  +7.21%  '
$ws.Range("D18").Value = '0.00001084'
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("D19").Value = '0.06563'
$ws.Range("E19").Value = '  +10.37%  '
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").Value = '17.54'
$ws.Range("E21").Value = '  +5.13%  '
$ws.Range("D22").Value = '6.225'
$ws.Range("E22").Value = '  +2.71%  '
$ws.Range("D23").Value = '28.287.31'
$ws.Range("E23").Value = '  +3.46%  '
$ws.Range("D24").Value = '11.69'
$ws.Range("E24").Value = '  +3.56%  '
$ws.Range("D25").Value = '2.175'
$ws.Range("E25").Value = '  -9.36%  '
$ws.Range("D26").Value = '20.79'
$ws.Range("E26").Value = '  +3.62%  '
$ws.Range("D27").Value = '156.28'
$ws.Range("E27").Value = '  +5.05%  '
$ws.Range("D28").Value = '2.034.59'
$ws.Range("E28").Value = '  +5.65%  '
$ws.Range("D29").Value = '2.322'
$ws.Range("E29").Value = '  -0.79%  '
$ws.Range("D30").Value = '128.32'
$ws.Range("E30").Value = '  +1.99%  '
$ws.Range("E31").Value = '  -1.54%  '
$ws.Range("D32").Value = '5.884'
$ws.Range("E32").Value = '  +5.33%  '
$ws.Range("D33").Value = '0.09245'
$ws.Range("E33").Value = '  +1.96%  '
$ws.Range("D34").Value = '3.657'
$ws.Range("E34").Value = '  +1.36%  '
$ws.Range("D35").Value = '13.09'
$ws.Range("E35").Value = '  +3.90%  '
$ws.Range("E36").Value = '  +4.88%  '
$ws.Range("D37").Value = '0.2186'
$ws.Range("E37").Value = '  +1.44%  '
$ws.Range("D38").Value = '5.188'
$ws.Range("E38").Value = '  +3.21%  '
$ws.Range("D39").Value = '0.06239'
$ws.Range("E39").Value = '  +2.36%  '
$ws.Range("D40").Value = '0.6582'
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.198'
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '8.175'
$ws.Range("E42").Value = '  +3.91%  '
$ws.Range("D43").Value = '0.9997'
$ws.Range("E43").Value = '  -0.17%  '
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").Value = '13.96'
$ws.Range("E45").Value = '  +3.58%  '
$ws.Range("D46").Value = '0.6094'
$ws.Range("E46").Value = '  +5.10%  '
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.10'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.08%  '
$ws.Range("D49").Value = '2.039'
$ws.Range("E49").Value = '  +5.49%  '
$ws.Range("E50").Value = '  +5.92%  '
$ws.Range("D51").Value = '0.06976'
$ws.Range("E51").Value = '  +2.44%  '
